# Adds two new demonstration tables after the last (bookmarked) table in the
# document:
#   1. a table whose rows start with empty grid positions (w:gridBefore)
#   2. a table whose rows end with empty grid positions (w:gridAfter)
# Both are inserted — preceded by a new empty paragraph — right after the
# document's final existing table and before the trailing empty paragraph /
# sectPr, mirroring the upstream python-docx `tbl-props.docx` test fixture.

$d = $word.ActiveDocument

# Locate the end of the last table currently in the document (the one that
# carries the "_GoBack" bookmark) so the new content lands right after it,
# regardless of any pre-existing offsets.
$lastTable = $d.Tables.Item($d.Tables.Count)
$insertAt = $lastTable.Range.End
$insertionRange = $d.Range($insertAt, $insertAt)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParagraph = "<w:p $wNs/>"

$tableGridBefore = "<w:tbl $wNs>" +
  "<w:tblGrid>" +
    "<w:gridCol w:w=`"500`"/>" +
    "<w:gridCol w:w=`"500`"/>" +
    "<w:gridCol w:w=`"500`"/>" +
    "<w:gridCol w:w=`"500`"/>" +
  "</w:tblGrid>" +
  "<w:tr>" +
    "<w:tc><w:p><w:r><w:t>a</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>b</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>c</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>d</w:t></w:r></w:p></w:tc>" +
  "</w:tr>" +
  "<w:tr>" +
    "<w:trPr><w:gridBefore w:val=`"1`"/></w:trPr>" +
    "<w:tc><w:p><w:r><w:t>e</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>f</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>g</w:t></w:r></w:p></w:tc>" +
  "</w:tr>" +
  "<w:tr>" +
    "<w:trPr><w:gridBefore w:val=`"2`"/></w:trPr>" +
    "<w:tc><w:p><w:r><w:t>h</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>i</w:t></w:r></w:p></w:tc>" +
  "</w:tr>" +
  "<w:tr>" +
    "<w:trPr><w:gridBefore w:val=`"3`"/></w:trPr>" +
    "<w:tc><w:p><w:r><w:t>j</w:t></w:r></w:p></w:tc>" +
  "</w:tr>" +
  "</w:tbl>"

$tableGridAfter = "<w:tbl $wNs>" +
  "<w:tblGrid>" +
    "<w:gridCol w:w=`"500`"/>" +
    "<w:gridCol w:w=`"500`"/>" +
    "<w:gridCol w:w=`"500`"/>" +
    "<w:gridCol w:w=`"500`"/>" +
  "</w:tblGrid>" +
  "<w:tr>" +
    "<w:tc><w:p><w:r><w:t>a</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>b</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>c</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>d</w:t></w:r></w:p></w:tc>" +
  "</w:tr>" +
  "<w:tr>" +
    "<w:trPr><w:gridAfter w:val=`"1`"/></w:trPr>" +
    "<w:tc><w:p><w:r><w:t>e</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>f</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>g</w:t></w:r></w:p></w:tc>" +
  "</w:tr>" +
  "<w:tr>" +
    "<w:trPr><w:gridAfter w:val=`"2`"/></w:trPr>" +
    "<w:tc><w:p><w:r><w:t>h</w:t></w:r></w:p></w:tc>" +
    "<w:tc><w:p><w:r><w:t>i</w:t></w:r></w:p></w:tc>" +
  "</w:tr>" +
  "</w:tbl>"

$payload = $newParagraph + $tableGridBefore + $tableGridAfter

$result = $insertionRange.InsertXML($payload)
